$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '30.297.50'
Set-TextValue "E2" '  +0.96%  '

Set-TextValue "D3" '1.921.89'
Set-TextValue "E3" '  +0.62%  '

Set-TextValue "D4" '1.000'
Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '0.8139'
Set-TextValue "E5" '  +2.95%  '

Set-TextValue "D6" '244.20'
Set-TextValue "E6" '  +0.98%  '

Set-TextValue "D7" '1.001'
Set-TextValue "E7" '  +0.01%  '

Set-TextValue "D8" '0.3266'

Set-TextValue "D9" '27.26'
Set-TextValue "E9" '  +3.59%  '

Set-TextValue "E10" '  +5.29%  '

Set-TextValue "D11" '0.7946'
Set-TextValue "E11" '  +7.03%  '

Set-TextValue "D12" '0.08112'

Set-TextValue "D13" '1.924.91'
Set-TextValue "E13" '  +0.85%  '

Set-TextValue "D14" '5.412'
Set-TextValue "E14" '  +4.36%  '

Set-TextValue "D15" '94.18'
Set-TextValue "E15" '  +1.21%  '

Set-TextValue "D16" '30.315.84'
Set-TextValue "E16" '  +1.03%  '

Set-TextValue "D17" '14.29'
Set-TextValue "E17" '  +2.29%  '

Set-TextValue "D18" '6.069'
Set-TextValue "E18" '  +3.49%  '

Set-TextValue "D19" '250.63'
Set-TextValue "E19" '  +2.04%  '

Set-TextValue "D20" '0.000007870'
Set-TextValue "E20" '  +1.74%  '

Set-TextValue "D21" '2.181.63'
Set-TextValue "E21" '  +1.24%  '

Set-TextValue "D23" '8.010'
Set-TextValue "E23" '  +17.22%  '

Set-TextValue "E24" '  +0.03%  '

Set-TextValue "D25" '0.1681'
Set-TextValue "E25" '  +21.18%  '

Set-TextValue "D26" '9.500'
Set-TextValue "E26" '  +2.94%  '

Set-TextValue "D27" '167.70'

Set-TextValue "E28" '  +0.79%  '

Set-TextValue "D29" '2.157'
Set-TextValue "E29" '  +6.15%  '

Set-TextValue "E30" '  +0.50%  '

Set-TextValue "D31" '1.552'
Set-TextValue "E31" '  +2.54%  '

Set-TextValue "D32" '4.353'
Set-TextValue "E32" '  +0.94%  '

Set-TextValue "D33" '0.05691'
Set-TextValue "E33" '  +3.30%  '

Set-TextValue "D34" '4.148'
Set-TextValue "E34" '  +1.62%  '

Set-TextValue "D35" '1.302'
Set-TextValue "E35" '  +3.54%  '

Set-TextValue "D36" '0.7474'
Set-TextValue "E36" '  +2.14%  '

Set-TextValue "D37" '1.000'
Set-TextValue "E37" '  +0.09%  '

Set-TextValue "D38" '2.725'
Set-TextValue "E38" '  +0.15%  '

Set-TextValue "D39" '0.01963'
Set-TextValue "E39" '  +2.04%  '

Set-TextValue "E40" '  +1.38%  '

Set-TextValue "E41" '  +2.22%  '

Set-TextValue "D42" '74.82'
Set-TextValue "E42" '  +3.53%  '

Set-TextValue "D43" '5.996'
Set-TextValue "E43" '  -2.08%  '

Set-TextValue "D44" '0.8558'
Set-TextValue "E44" '  +2.21%  '

Set-TextValue "D45" '1.929'
Set-TextValue "E45" '  +2.98%  '

Set-TextValue "D46" '1.001'
Set-TextValue "E46" '  +0.02%  '

Set-TextValue "D47" '1.038.26'
Set-TextValue "E47" '  +5.23%  '

Set-TextValue "D48" '103.04'
Set-TextValue "E48" '  +2.70%  '

Set-TextValue "E49" '  +11.76%  '

Set-TextValue "B50" 'Aptos'
Set-TextValue "C50" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D50" '7.659'
Set-TextValue "E50" '  +1.50%  '

Set-TextValue "B51" 'EnergySwap'
Set-TextValue "C51" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D51" '9.931'
Set-TextValue "E51" '  +1.31%  '
